$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated stress/deflection values in column B (rows 2-11)
$values = @{
    3  = 1.1539999999999999
    4  = 1.256
    5  = 1.5
    6  = 1.605
    7  = 1.7050000000000001
    8  = 1.804
    9  = 1.9039999999999999
    10 = 1.9530000000000001
    11 = 2
}

foreach ($row in $values.Keys) {
    $ws.Range("B$row").Value = $values[$row]
}

# Add new (currently empty) column D cells, rows 3-11, formatted with the
# same numeric format used elsewhere in the table (0.000) - this creates
# the new cellXfs style used by these cells.
for ($row = 3; $row -le 11; $row++) {
    $ws.Range("D$row").NumberFormat = "0.000"
}

# Move the selection like the authored workbook (selection now sits on D22)
$ws.Range("D22").Select()
